# Apply the changes described by the diff to the "Task detail" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task detail")

# 1) Update the text in C11:
#    - move "- ok" on line 3 (after the parenthetical note instead of before it)
#    - add "- ok" to the end of line 6
$newText = "1/ Khi chuyển qua từng user chat, làm hiệu ứng chuyển bên khung chat - ok`n" +
           "2/ Mang link socket qua file config - ok`n" +
           "3/ Disconnect khi unmount vuejs component (tự động hủy socket khi change component) - ok`n" +
           "4/ Check duplicate email - ok`n" +
           "5/ Sửa lại logic script chổ resize, tối ưu xóa biến thừa - ok`n" +
           "6/ Tách input và button send ra thì 1 component riêng - ok`n" +
           " 7/ Room chat khi load lại, tin nhắn do mình gửi lại nằm bên trái - ok`n"

$ws.Range("C11").Value = $newText

# 2) Move the active selection from C12 to C11
$ws.Range("C11").Select()

# 3) Adjust row heights
$ws.Rows.Item(6).RowHeight = 75
$ws.Rows.Item(10).RowHeight = 150
